$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.434.11'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').Value = '3.326.76'
$ws.Range('E3').Value = '  -4.15%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '574.93'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '177.90'
$ws.Range('E6').Value = '  +2.56%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.614'
$ws.Range('E7').Value = '  +2.87%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '3.322.31'
$ws.Range('E9').Value = '  -4.19%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.129'
$ws.Range('E10').Value = '  -1.16%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.87'
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.412'
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').Value = '3.904.75'
$ws.Range('E13').Value = '  -4.12%  '
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '28.68'
$ws.Range('E15').Value = '  -4.11%  '
$ws.Range('D16').Value = '65.532.88'
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000169'
$ws.Range('E17').Value = '  -1.43%  '
$ws.Range('D18').Value = '3.328.75'
$ws.Range('E18').Value = '  -4.04%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.72'
$ws.Range('E19').Value = '  -3.47%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.41'
$ws.Range('E20').Value = '  -3.34%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '361.99'
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.42'
$ws.Range('E22').Value = '  -3.88%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '71.24'
$ws.Range('E24').Value = '  -1.46%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.519'
$ws.Range('E25').Value = '  -2.76%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0000122'
$ws.Range('E26').Value = '  -1.54%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.57'
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.61'
$ws.Range('E31').Value = '  -2.36%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '22.81'
$ws.Range('E33').Value = '  -4.34%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.83'
$ws.Range('E34').Value = '  -4.10%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.22'
$ws.Range('E35').Value = '  -4.90%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.49'
$ws.Range('E36').Value = '  -2.67%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '160.56'
$ws.Range('E37').Value = '  +0.89%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.844'
$ws.Range('E38').Value = '  -5.00%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '27.35'
$ws.Range('E39').Value = '  -6.03%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.75'
$ws.Range('E40').Value = '  -0.92%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.54'
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('D42').Value = '2.693.23'
$ws.Range('E42').Value = '  -3.99%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '6.22'
$ws.Range('E43').Value = '  -3.22%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.28'
$ws.Range('E44').Value = '  -3.71%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0666'
$ws.Range('E45').Value = '  -2.14%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '334.69'
$ws.Range('E46').Value = '  +7.91%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '39.63'
$ws.Range('E47').Value = '  -0.92%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '24.12'
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0278'
$ws.Range('E49').Value = '  -3.50%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.103'
$ws.Range('E50').Value = '  +2.76%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.963'
$ws.Range('E51').Value = '  -1.13%  '
